$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 with new values
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 206

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 205

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 201

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 18

# Remove the old row 6 entirely (A6/B6), shrinking the used range to A1:B5
$ws.Range("A6:B6").Delete()
